$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 256.1875
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H41").Value = 788.8823
$ws.Range("I41").Value = 606.75
$ws.Range("J41").Value = 950.7778
$ws.Range("K41").Value = 606.75
$ws.Range("L41").Value = 950.7778
$ws.Range("M41").Value = -166.75
$ws.Range("N41").Value = -1830.7778
$ws.Range("H51").Value = 4079.1667
$ws.Range("J51").Value = 4100
$ws.Range("L51").Value = 4100
$ws.Range("N51").Value = -5068
$ws.Range("H86").Value = 4364.4
$ws.Range("I86").Value = 4215.3335
$ws.Range("K86").Value = 4215.3335
$ws.Range("M86").Value = -3092.3335
$ws.Range("H89").Value = 4364.4
$ws.Range("I89").Value = 4215.3335
$ws.Range("K89").Value = 21076.6675
$ws.Range("M89").Value = -15460.6675
$ws.Range("H92").Value = 559.93335
$ws.Range("I92").Value = 483.5
$ws.Range("J92").Value = 865.6667
$ws.Range("K92").Value = 483.5
$ws.Range("L92").Value = 865.6667
$ws.Range("M92").Value = 764.5
$ws.Range("N92").Value = -3361.6667
$ws.Range("H100").Value = 4049.5
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 5099
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 5099
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -6181
$ws.Range("H106").Value = 2675
$ws.Range("I106").Value = 2650
$ws.Range("K106").Value = 2650
$ws.Range("M106").Value = -2019
$ws.Range("H108").Value = 99953.71000000001
$ws.Range("J108").Value = 99953.71000000001
$ws.Range("L108").Value = 99953.71000000001
$ws.Range("N108").Value = -107633.71
$ws.Range("H116").Value = 8075.1177
$ws.Range("I116").Value = 7600.25
$ws.Range("J116").Value = 8497.223
$ws.Range("K116").Value = 7600.25
$ws.Range("L116").Value = 8497.223
$ws.Range("M116").Value = -4158.25
$ws.Range("N116").Value = -15381.223
$ws.Range("H117").Value = 99969.86
$ws.Range("J117").Value = 99969.86
$ws.Range("L117").Value = 99969.86
$ws.Range("N117").Value = -109147.86
$ws.Range("H137").Value = 428877.88
$ws.Range("I137").Value = 2321.9443
$ws.Range("J137").Value = 908753.3
$ws.Range("K137").Value = 6965.8329
$ws.Range("L137").Value = 2726259.9
$ws.Range("M137").Value = -4415.8329
$ws.Range("N137").Value = -2731359.9

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8335928
$ws.Range("I45").Value = 2738.875
$ws.Range("J45").Value = 17859572
$ws.Range("K45").Value = 2738.875
$ws.Range("L45").Value = 17859572
$ws.Range("M45").Value = -2361.875
$ws.Range("N45").Value = -17860326
$ws.Range("H61").Value = 39984.77
$ws.Range("I61").Value = 1277.5454
$ws.Range("K61").Value = 1277.5454
$ws.Range("M61").Value = -1065.5454
$ws.Range("H122").Value = 2738.1516
$ws.Range("I122").Value = 2173.963
$ws.Range("J122").Value = 5277
$ws.Range("K122").Value = 6521.889000000001
$ws.Range("L122").Value = 15831
$ws.Range("M122").Value = -4071.889000000001
$ws.Range("N122").Value = -20731
$ws.Range("H132").Value = 2001.7858
$ws.Range("J132").Value = 2326.5
$ws.Range("L132").Value = 6979.5
$ws.Range("N132").Value = -12039.5
$ws.Range("H135").Value = 59497.332
$ws.Range("J135").Value = 59497.332
$ws.Range("L135").Value = 59497.332
$ws.Range("N135").Value = -69637.33199999999
$ws.Range("H136").Value = 39984.77
$ws.Range("I136").Value = 1277.5454
$ws.Range("K136").Value = 3832.6362
$ws.Range("M136").Value = -1282.6362
$ws.Range("H139").Value = 78994
$ws.Range("J139").Value = 78994
$ws.Range("L139").Value = 78994
$ws.Range("N139").Value = -89274
$ws.Range("H141").Value = 91831.336
$ws.Range("J141").Value = 80997
$ws.Range("L141").Value = 80997
$ws.Range("N141").Value = -91357

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5584.5386
$ws.Range("J86").Value = 10995.5
$ws.Range("L86").Value = 10995.5
$ws.Range("N86").Value = -13241.5
$ws.Range("H89").Value = 5584.5386
$ws.Range("J89").Value = 10995.5
$ws.Range("L89").Value = 54977.5
$ws.Range("N89").Value = -66209.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2684.7273
$ws.Range("I31").Value = 1971.4706
$ws.Range("K31").Value = 1971.4706
$ws.Range("M31").Value = -1676.4706
$ws.Range("H34").Value = 2684.7273
$ws.Range("I34").Value = 1971.4706
$ws.Range("K34").Value = 1971.4706
$ws.Range("M34").Value = -1769.4706
$ws.Range("H58").Value = 1239.697
$ws.Range("I58").Value = 1189.1538
$ws.Range("K58").Value = 1189.1538
$ws.Range("M58").Value = -986.1538
$ws.Range("H63").Value = 55000
$ws.Range("J63").Value = 55000
$ws.Range("L63").Value = 55000
$ws.Range("N63").Value = -56372
$ws.Range("H66").Value = 55000
$ws.Range("J66").Value = 55000
$ws.Range("L66").Value = 165000
$ws.Range("N66").Value = -171864
$ws.Range("H136").Value = 1239.697
$ws.Range("I136").Value = 1189.1538
$ws.Range("K136").Value = 3567.4614
$ws.Range("M136").Value = -1017.4614

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 13.153846
$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 13.285714
$ws.Range("K2").Value = 78
$ws.Range("L2").Value = 79.71428400000001
$ws.Range("M2").Value = 35
$ws.Range("N2").Value = -305.714284
$ws.Range("H7").Value = 6001.5713
$ws.Range("J7").Value = 10467.5
$ws.Range("L7").Value = 31402.5
$ws.Range("N7").Value = -31626.5
$ws.Range("H23").Value = 125049.75
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 142910.58
$ws.Range("K23").Value = 72
$ws.Range("L23").Value = 428731.74
$ws.Range("M23").Value = 163
$ws.Range("N23").Value = -429201.74
$ws.Range("H34").Value = 2746.4167
$ws.Range("I34").Value = 423.33334
$ws.Range("J34").Value = 3520.7778
$ws.Range("K34").Value = 1270.00002
$ws.Range("L34").Value = 10562.3334
$ws.Range("M34").Value = -1186.00002
$ws.Range("N34").Value = -10730.3334
$ws.Range("H38").Value = 607.55554
$ws.Range("I38").Value = 39
$ws.Range("J38").Value = 1062.4
$ws.Range("K38").Value = 117
$ws.Range("L38").Value = 3187.2
$ws.Range("M38").Value = 230
$ws.Range("N38").Value = -3881.2
$ws.Range("H39").Value = 1621.35
$ws.Range("J39").Value = 1834
$ws.Range("L39").Value = 5502
$ws.Range("N39").Value = -6090
$ws.Range("H107").Value = 1466.3334
$ws.Range("J107").Value = 1249.5
$ws.Range("L107").Value = 3748.5
$ws.Range("N107").Value = -7588.5
$ws.Range("H133").Value = 13419
$ws.Range("I133").Value = 5114.222
$ws.Range("K133").Value = 15342.666
$ws.Range("M133").Value = -10282.666

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 786
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 1014.6667
$ws.Range("K13").Value = 100
$ws.Range("L13").Value = 1014.6667
$ws.Range("M13").Value = 39
$ws.Range("N13").Value = -1292.6667
$ws.Range("H70").Value = 56514.9
$ws.Range("I70").Value = 62017.668
$ws.Range("K70").Value = 62017.668
$ws.Range("M70").Value = -61747.668
$ws.Range("H73").Value = 56514.9
$ws.Range("I73").Value = 62017.668
$ws.Range("K73").Value = 62017.668
$ws.Range("M73").Value = -61081.668
$ws.Range("H102").Value = 1867.0416
$ws.Range("J102").Value = 1786
$ws.Range("L102").Value = 1786
$ws.Range("N102").Value = -5030
$ws.Range("H122").Value = 347968.6
$ws.Range("I122").Value = 437312.25
$ws.Range("K122").Value = 1311936.75
$ws.Range("M122").Value = -1309486.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 56245.43
$ws.Range("J133").Value = 56245.43
$ws.Range("L133").Value = 56245.43
$ws.Range("N133").Value = -61305.43
$ws.Range("H135").Value = 78065.336
$ws.Range("J135").Value = 78065.336
$ws.Range("L135").Value = 78065.336
$ws.Range("N135").Value = -88205.336
$ws.Range("H138").Value = 98877.60000000001
$ws.Range("J138").Value = 98877.60000000001
$ws.Range("L138").Value = 98877.60000000001
$ws.Range("N138").Value = -109157.6

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 15000444
$ws.Range("I2").Value = 15000444
$ws.Range("K2").Value = 15000444
$ws.Range("M2").Value = -15000332
$ws.Range("H46").Value = 96002
$ws.Range("J46").Value = 96002
$ws.Range("L46").Value = 96002
$ws.Range("N46").Value = -96464
$ws.Range("H134").Value = 96002
$ws.Range("J134").Value = 96002
$ws.Range("L134").Value = 288006
$ws.Range("N134").Value = -293076
